# Moves a set of nine shapes + one diagram (graphicFrame) on slide 2 by the
# same fixed offset (35169 EMU right, 34728 EMU down), as if the author had
# selected them together and nudged/dragged them as one unit.
#
# Shape.Left / Shape.Top in this object model are backed by a 32-bit
# (Single) float, so setting each shape's Left/Top independently with a
# naive point value cannot reliably reproduce an exact target EMU offset
# (float32 has ~7 significant digits; at these coordinate magnitudes one
# ULP can span more than 1 EMU). To dodge that, every affected shape is
# temporarily collected into one group; only the *group's* Left/Top is
# set (a single float32 rounding event), and PowerPoint's group->child
# transform propagates that as an exact integer EMU delta to every child
# (off stays an untouched integer; only the group off/chOff pair is
# float-rounded once). The group is then dissolved again and original
# z-order is restored so the slide XML only differs in the ten <a:off>
# values, exactly like the source diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeIndexById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        if ($slide.Shapes.Item($i).Id -eq $id) {
            return $i
        }
    }
    return -1
}

# IDs (p:cNvPr @id) of every shape that moves, and the slide's full
# original front-to-back order (by Id) so it can be restored afterwards.
$moveIds     = @(11, 6, 17, 24, 25, 26, 27, 28, 4, 23)
$originalOrder = @(11, 6, 7, 13, 17, 18, 24, 25, 26, 27, 28, 29, 30, 31, 4, 2, 23)

# Build the COM index array Shapes.Range() needs, resolving each Id to its
# current position right now (before any grouping moves things around).
$moveIdxs = @()
foreach ($id in $moveIds) {
    $moveIdxs += (Get-ShapeIndexById $s $id)
}

$range = $s.Shapes.Range($moveIdxs)
$grp = $range.Group()

# Target group bounding-box origin, precomputed so that casting to the
# COM layer's Single and flooring(*12700) lands exactly on the OOXML EMU
# values from the target diff (1612533, 3493177) -- i.e. the original
# bounding-box origin (1577364, 3458449) shifted by (+35169, +34728) EMU.
$grp.Left = 126.97110748291016
$grp.Top = 275.0533142089844

$grp.Ungroup() | Out-Null

# Restore original z-order: sending shapes to back in reverse target
# order rebuilds the exact original front-to-back sequence.
for ($k = $originalOrder.Count; $k -ge 1; $k--) {
    $id = $originalOrder[$k - 1]
    $idx = Get-ShapeIndexById $s $id
    $s.Shapes.Item($idx).ZOrder(1) | Out-Null   # msoSendToBack
}
